$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 8435
$wsExhibit.Range("F5").Value = 6145
$wsExhibit.Range("F11").Value = 1116

# Sheet "全部类型" (all types) - same rows duplicated, update matching cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8435
$wsAll.Range("F5").Value = 6145
$wsAll.Range("F15").Value = 1116
